# Update odds values on Sheet1 as described in the commit diff.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Row 4
$ws.Range("M4").Value = 1.1
$ws.Range("O4").Value = 1.44
$ws.Range("P4").Value = 2.63

# Row 5
$ws.Range("V5").Value = 1.73

# Row 6
$ws.Range("U6").Value = 1.67

# Row 7
$ws.Range("L7").Value = 4.5
$ws.Range("S7").Value = 1.5
$ws.Range("T7").Value = 2.5
$ws.Range("V7").Value = 1.73
$ws.Range("AT7").Value = 2.5
